$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"

$ws.Range("D2").Value = "0"

$ws.Range("C3").Value = "0"
$ws.Range("D3").Value = "1"

$ws.Range("C4").Value = "3"
$ws.Range("D4").Value = "4"
